# Apply the "updated summary charts and summary reports including
# comments from Prof. Erhardt" edit to the Stockton, CA Metro Area-Bus
# FAC summary workbook.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# ---------------------------------------------------------------------
# Sheet1: text / year updates
# ---------------------------------------------------------------------
$ws1.Range("C1").Value = "2012"
$ws1.Range("E7").Value = "2012"

# ---------------------------------------------------------------------
# Sheet1: updated factor data (rows 8-18), formulas switched from
# "*100/" (percentage-as-number) to a plain ratio now displayed with a
# percentage number format, and refreshed Ridership Effect values.
# ---------------------------------------------------------------------

# Row 8 - Vehicle Revenue Miles
$ws1.Range("E8").Value = 2816957
$ws1.Range("F8").Value = 2471911
$ws1.Range("G8").Formula = "=IFERROR((F8-E8)/E8,0)"
$ws1.Range("H8").Value = -792329.41021
$ws1.Range("I8").Formula = "=IFERROR(H8/$E$21,0)"

# Row 9 - Average Fare (2018$)
$ws1.Range("E9").Value = 1.205051927
$ws1.Range("F9").Value = 0.894024738
$ws1.Range("G9").Formula = "=IFERROR((F9-E9)/E9,0)"
$ws1.Range("H9").Value = -270797.89291
$ws1.Range("I9").Formula = "=IFERROR(H9/$E$21,0)"

# Row 10 - Population + Employment
$ws1.Range("E10").Value = 969671
$ws1.Range("F10").Value = 1064707.73
$ws1.Range("G10").Formula = "=IFERROR((F10-E10)/E10,0)"
$ws1.Range("H10").Value = 169539.10345
$ws1.Range("I10").Formula = "=IFERROR(H10/$E$21,0)"

# Row 11 - % of Population in Transit Supportive Density
$ws1.Range("E11").Value = 47.04104114
$ws1.Range("F11").Value = 45.74463935
$ws1.Range("G11").Formula = "=IFERROR((F11-E11)/E11,0)"
$ws1.Range("H11").Value = -43320.704426
$ws1.Range("I11").Formula = "=IFERROR(H11/$E$21,0)"

# Row 12 - Average Gas Price (2018$)
$ws1.Range("E12").Value = 4.3491
$ws1.Range("F12").Value = 3.4
$ws1.Range("G12").Formula = "=IFERROR((F12-E12)/E12,0)"
$ws1.Range("H12").Value = -163123.17908
$ws1.Range("I12").Formula = "=IFERROR(H12/$E$21,0)"

# Row 13 - Median Per Capita (2018$)
$ws1.Range("E13").Value = 24917.4
$ws1.Range("F13").Value = 30789.5
$ws1.Range("G13").Formula = "=IFERROR((F13-E13)/E13,0)"
$ws1.Range("H13").Value = -231149.06125
$ws1.Range("I13").Formula = "=IFERROR(H13/$E$21,0)"

# Row 14 - % of Households with 0 Vehicles
$ws1.Range("E14").Value = 7.21
$ws1.Range("F14").Value = 3.52
$ws1.Range("G14").Formula = "=IFERROR((F14-E14)/E14,0)"
$ws1.Range("H14").Value = -115289.020014
$ws1.Range("I14").Formula = "=IFERROR(H14/$E$21,0)"

# Row 15 - % Working at Home
$ws1.Range("E15").Value = 3.5
$ws1.Range("F15").Value = 3.4
$ws1.Range("G15").Formula = "=IFERROR((F15-E15)/E15,0)"
$ws1.Range("H15").Value = 6862.407423
$ws1.Range("I15").Formula = "=IFERROR(H15/$E$21,0)"

# Row 16 - Years Since Ride-hail Start (E16/F16 remain blank)
$ws1.Range("G16").Formula = "=IFERROR((F16-E16)/E16,0)"
$ws1.Range("H16").Value = -566023.69655
$ws1.Range("I16").Formula = "=IFERROR(H16/$E$21,0)"

# Row 17 - Bike Share
$ws1.Range("E17").Value = 0
$ws1.Range("F17").Value = 0
$ws1.Range("G17").Formula = "=IFERROR((F17-E17)/E17,0)"
$ws1.Range("H17").Value = 0
$ws1.Range("I17").Formula = "=IFERROR(H17/$E$21,0)"

# Row 18 - Electric Scooters
$ws1.Range("E18").Value = 0
$ws1.Range("F18").Value = 0
$ws1.Range("G18").Formula = "=IFERROR((F18-E18)/E18,0)"
$ws1.Range("H18").Value = 0
$ws1.Range("I18").Formula = "=IFERROR(H18/$E$21,0)"

# Row 19 - New Reporters (now carries an explicit 0 Ridership Effect)
$ws1.Range("G19").Formula = "=IFERROR((F19-E19)/E19,0)"
$ws1.Range("H19").Value = 0
$ws1.Range("I19").Formula = "=IFERROR(H19/$E$21,0)"

# Row 20 - Total Modeled Ridership
$ws1.Range("E20").Value = 4089204.814
$ws1.Range("F20").Value = 3198847.89
$ws1.Range("G20").Formula = "=IFERROR((F20-E20)/E20,0)"
$ws1.Range("I20").Formula = "=G20"

# Row 21 - Total Observed Ridership
$ws1.Range("E21").Value = 4257409
$ws1.Range("F21").Value = 3510022
$ws1.Range("G21").Formula = "=IFERROR((F21-E21)/E21,0)"
$ws1.Range("I21").Formula = "=G21"

# ---------------------------------------------------------------------
# Number-format refresh: ratios now render as percentages, and the
# average-value / ridership-effect columns render with 2 decimals.
# ---------------------------------------------------------------------
$ws1.Range("E8:F18").NumberFormat = "#,##0.00"
$ws1.Range("H8:H18").NumberFormat = "#,##0.00"
$ws1.Range("G8:G21").NumberFormat = "0.00%"
$ws1.Range("I8:I19").NumberFormat = "0.00%"

$ws1.Range("E19:F19").NumberFormat = "#,##0.00"
$ws1.Range("H19").NumberFormat = "#,##0.00"
$ws1.Range("E21:F21").NumberFormat = "#,##0.00"
$ws1.Range("H21").NumberFormat = "#,##0.00"
$ws1.Range("I21").NumberFormat = "0.00%"

$ws1.Range("E20:F20").NumberFormat = "#,##0.00"
$ws1.Range("H20").NumberFormat = "#,##0.00"
$ws1.Range("G20").NumberFormat = "0.00%"
$ws1.Range("I20").NumberFormat = "0.00%"

# ---------------------------------------------------------------------
# Sheet1 view: scroll back to the top and move the selection.
# ---------------------------------------------------------------------
$ws1.Activate()
$ws1.Range("H21").Select()

Write-Host "Edit applied"
